$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.142.41"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").Value = "3.527.06"
$ws.Range("E3").Value = "  +3.63%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").Value = "593.01"
$ws.Range("E5").Value = "  +0.92%  "

# Row 6
$ws.Range("D6").Value = "138.62"
$ws.Range("E6").Value = "  +4.89%  "

# Row 7
$ws.Range("D7").Value = "3.525.91"
$ws.Range("E7").Value = "  +3.62%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  +2.68%  "

# Row 10
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  +5.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.79%  "

# Row 12
$ws.Range("D12").Value = "0.388"
$ws.Range("E12").Value = "  +5.08%  "

# Row 13
$ws.Range("D13").Value = "4.114.40"
$ws.Range("E13").Value = "  +3.27%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000186"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.92%  "

# Row 15
$ws.Range("D15").Value = "26.96"
$ws.Range("E15").Value = "  +4.04%  "

# Row 16
$ws.Range("D16").Value = "3.524.45"
$ws.Range("E16").Value = "  +3.38%  "

# Row 17
$ws.Range("E17").Value = "  +1.26%  "

# Row 18
$ws.Range("D18").Value = "64.968.43"
$ws.Range("E18").Value = "  +0.00%  "

# Row 19
$ws.Range("D19").Value = "10.12"
$ws.Range("E19").Value = "  +1.23%  "

# Row 20
$ws.Range("D20").Value = "5.85"
$ws.Range("E20").Value = "  +4.17%  "

# Row 21
$ws.Range("D21").Value = "14.21"
$ws.Range("E21").Value = "  +5.90%  "

# Row 22
$ws.Range("D22").Value = "394.34"
$ws.Range("E22").Value = "  +1.64%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.570"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.51%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.00%  "

# Row 25
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.656.02"
$ws.Range("E25").Value = "  +3.08%  "

# Row 26
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000113"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.80%  "

# Row 28
$ws.Range("D28").Value = "7.68"
$ws.Range("E28").Value = "  +10.21%  "

# Row 29
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.09%  "

# Row 30
$ws.Range("D30").Value = "2.25"
$ws.Range("E30").Value = "  +0.93%  "

# Row 31
$ws.Range("D31").Value = "8.27"
$ws.Range("E31").Value = "  +3.79%  "

# Row 32
$ws.Range("D32").Value = "3.536.27"
$ws.Range("E32").Value = "  +3.74%  "

# Row 33
$ws.Range("E33").Value = "  -0.02%  "

# Row 34
$ws.Range("D34").Value = "23.71"
$ws.Range("E34").Value = "  +7.13%  "

# Row 35
$ws.Range("E35").Value = "  +2.59%  "

# Row 36
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "1.22"
$ws.Range("E36").Value = "  +1.37%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "169.79"
$ws.Range("E37").Value = "  -1.95%  "

# Row 38
$ws.Range("D38").Value = "6.94"
$ws.Range("E38").Value = "  +3.81%  "

# Row 39
$ws.Range("D39").Value = "1.53"
$ws.Range("E39").Value = "  +3.18%  "

# Row 40
$ws.Range("D40").Value = "4.83"
$ws.Range("E40").Value = "  +3.28%  "

# Row 41
$ws.Range("D41").Value = "0.0794"
$ws.Range("E41").Value = "  +6.53%  "

# Row 42
$ws.Range("D42").Value = "0.819"
$ws.Range("E42").Value = "  +1.96%  "

# Row 43
$ws.Range("D43").Value = "26.69"
$ws.Range("E43").Value = "  +26.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.75%  "

# Row 45
$ws.Range("E45").Value = "  -0.23%  "

# Row 46
$ws.Range("D46").Value = "4.41"
$ws.Range("E46").Value = "  +1.23%  "

# Row 47
$ws.Range("E47").Value = "  +13.68%  "

# Row 48
$ws.Range("E48").Value = "  +5.09%  "

# Row 49
$ws.Range("D49").Value = "6.77"
$ws.Range("E49").Value = "  +5.20%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.294.48"
$ws.Range("E50").Value = "  +5.32%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.59%  "
